# Updates the crypto price/volume table on Sheet1 for the Thu Mar 30 2023
# GitHub Actions data refresh. Column D ("Price") and column E ("Volume(1h)")
# are plain text cells in the source sheet (e.g. "28.565.71", "  +0.51%  "),
# so for any Column D value that looks like a normal decimal number we
# prefix it with a leading apostrophe to force Excel to keep storing it as
# text instead of silently reinterpreting it as a Number (which would also
# strip meaningful trailing zeros, e.g. "0.06530" -> "0.0653"). Values that
# already contain two dots (e.g. "28.565.71") are never parsed as numbers by
# Excel, so they do not need the apostrophe.
#
# Rows 33/34 and 40/41 also swapped coin identity (Filecoin <-> HuobiToken,
# TheSandbox <-> Aptos), so their Coin name (B) and Link (C) columns are
# rewritten too, in addition to Price/Volume.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    $ws.Range("D2").Value = '28.565.71'
    $ws.Range("E2").Value = '  +0.51%  '
    $ws.Range("D3").Value = '1.800.73'
    $ws.Range("E3").Value = '  -0.49%  '
    $ws.Range("D4").Value = '''1.002'
    $ws.Range("E4").Value = '  +0.10%  '
    $ws.Range("D5").Value = '''317.41'
    $ws.Range("E5").Value = '  +0.17%  '
    $ws.Range("D7").Value = '''0.5418'
    $ws.Range("E7").Value = '  -1.75%  '
    $ws.Range("D8").Value = '''0.3782'
    $ws.Range("E8").Value = '  -1.94%  '
    $ws.Range("D9").Value = '''0.07494'
    $ws.Range("E9").Value = '  -1.37%  '
    $ws.Range("D10").Value = '''42.18'
    $ws.Range("E10").Value = '  -1.87%  '
    $ws.Range("D11").Value = '''1.111'
    $ws.Range("E11").Value = '  -1.65%  '
    $ws.Range("D12").Value = '''1.002'
    $ws.Range("E12").Value = '  +0.15%  '
    $ws.Range("D13").Value = '''20.63'
    $ws.Range("E13").Value = '  -2.34%  '
    $ws.Range("D14").Value = '''6.155'
    $ws.Range("E14").Value = '  -0.94%  '
    $ws.Range("D15").Value = '''7.329'
    $ws.Range("E15").Value = '  +0.01%  '
    $ws.Range("D16").Value = '1.794.51'
    $ws.Range("E16").Value = '  -0.33%  '
    $ws.Range("D17").Value = '''89.99'
    $ws.Range("E17").Value = '  -1.24%  '
    $ws.Range("D18").Value = '''0.00001066'
    $ws.Range("E18").Value = '  -0.84%  '
    $ws.Range("D19").Value = '''0.06530'
    $ws.Range("E19").Value = '  +1.01%  '
    $ws.Range("D20").Value = '''17.54'
    $ws.Range("E20").Value = '  +1.65%  '
    $ws.Range("E21").Value = '  +0.10%  '
    $ws.Range("D22").Value = '''5.941'
    $ws.Range("E22").Value = '  -0.82%  '
    $ws.Range("D23").Value = '28.557.68'
    $ws.Range("E23").Value = '  +0.43%  '
    $ws.Range("D24").Value = '''11.13'
    $ws.Range("E24").Value = '  -1.48%  '
    $ws.Range("D25").Value = '''2.090'
    $ws.Range("E25").Value = '  -0.92%  '
    $ws.Range("D26").Value = '''160.20'
    $ws.Range("E26").Value = '  +2.59%  '
    $ws.Range("D27").Value = '''20.48'
    $ws.Range("E27").Value = '  -1.14%  '
    $ws.Range("D28").Value = '2.003.24'
    $ws.Range("E28").Value = '  -0.47%  '
    $ws.Range("D29").Value = '''2.337'
    $ws.Range("E29").Value = '  -4.01%  '
    $ws.Range("D30").Value = '''122.97'
    $ws.Range("E30").Value = '  -0.81%  '
    $ws.Range("D31").Value = '''1.124'
    $ws.Range("E31").Value = '  -4.16%  '
    $ws.Range("D32").Value = '''0.1059'
    $ws.Range("E32").Value = '  +2.36%  '
    $ws.Range("B33").Value = 'HuobiToken'
    $ws.Range("C33").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
    $ws.Range("D33").Value = '''3.688'
    $ws.Range("E33").Value = '  +1.15%  '
    $ws.Range("B34").Value = 'Filecoin'
    $ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    $ws.Range("D34").Value = '''5.627'
    $ws.Range("E34").Value = '  -2.18%  '
    $ws.Range("D35").Value = '''0.06559'
    $ws.Range("E35").Value = '  +4.88%  '
    $ws.Range("D36").Value = '''0.2251'
    $ws.Range("E36").Value = '  -2.53%  '
    $ws.Range("D37").Value = '''0.02300'
    $ws.Range("E37").Value = '  -1.30%  '
    $ws.Range("D38").Value = '''8.636'
    $ws.Range("E38").Value = '  -3.33%  '
    $ws.Range("D39").Value = '''5.030'
    $ws.Range("E39").Value = '  -0.13%  '
    $ws.Range("B40").Value = 'Aptos'
    $ws.Range("C40").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
    $ws.Range("D40").Value = '''11.23'
    $ws.Range("E40").Value = '  -3.16%  '
    $ws.Range("B41").Value = 'TheSandbox'
    $ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
    $ws.Range("D41").Value = '''0.6215'
    $ws.Range("E41").Value = '  -2.80%  '
    $ws.Range("E42").Value = '  +2.33%  '
    $ws.Range("D43").Value = '''1.452'
    $ws.Range("E43").Value = '  +4.64%  '
    $ws.Range("E44").Value = '  +0.07%  '
    $ws.Range("D45").Value = '''13.36'
    $ws.Range("E45").Value = '  -0.48%  '
    $ws.Range("D46").Value = '''3.693'
    $ws.Range("E46").Value = '  +0.13%  '
    $ws.Range("D47").Value = '''0.5842'
    $ws.Range("E47").Value = '  -2.56%  '
    $ws.Range("D48").Value = '''127.29'
    $ws.Range("E48").Value = '  +2.60%  '
    $ws.Range("D49").Value = '''1.957'
    $ws.Range("E49").Value = '  -0.72%  '
    $ws.Range("D50").Value = '''1.199'
    $ws.Range("E50").Value = '  +4.50%  '
    $ws.Range("D51").Value = '''0.06899'
    $ws.Range("E51").Value = '  -0.47%  '
